$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Column A (rows 3-6) held "Pontos" numbers 1-4 stored as text; convert to real numbers.
#    (Row 7 already stored its "5" as a real number.)
$ws.Range("A3:A6").NumberFormat = "General"
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4

# 2. Add unit suffixes to the column headers in row 2.
$ws.Range("B2").Value = "Latitude (º)"
$ws.Range("C2").Value = "Longitude (º)"
$ws.Range("D2").Value = "Altura (m)"

# 3. Update the saved selection to the full data block used for the new slide.
$ws.Range("A1:G7").Select()

# 4. Configure the page for printing (new slide export used A4/portrait).
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
